# Insert a new detail row for September (row 32) into the "2024" sheet,
# pushing all subsequent rows down by one. The new row records a
# "spotify" transaction dated 2024-09-06 16:34:27 in the
# September_Details / September_Date columns (R/S).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Shift row 32 (and everything below it) down by one row.
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new September entry.
$ws.Range("R32").Value = "spotify"
$ws.Range("S32").Value = "2024-09-06 16:34:27"
